# Scenario_Component_Behavior.xlsx — "make the model feasible" tweak.
#
# Per the commit: shrink two over-large annual-demand inputs on the
# OperationScenario_Behavior sheet (row 2) so the optimisation model stops
# being infeasible:
#   I2  hot_water_demand_annual           2,500,000 -> 2,700   (Wh)
#   L2  appliance_electricity_demand_annual 4,000,000 -> 1,000 (Wh)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 2700
$ws.Range("L2").Value = 1000

# Author's saved cursor/selection moved to I3 (and the window scrolled right
# so column H is the first visible column) when they made the edit above.
$ws.Range("I3").Select()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
